# Added a new pattern to underground levels that is based upon a pattern
# from underwater level types.
#
# Net effect: on Sheet1, the FDS row's "Official Name" cell (B6) is updated
# from "Super Mario Bros (J).FDS" to "Super Mario Bros. (J).fds", and the
# active selection moves to B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("B6").Value = "Super Mario Bros. (J).fds"

$ws.Range("B6").Select()
